$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "25.841.23"
$ws.Range("E2").Value = "  -0.72%  "
Set-TextValue $ws.Range("D3") "1.630.04"
$ws.Range("E3").Value = "  -0.81%  "
$ws.Range("E4").Value = "  +0.06%  "
Set-TextValue $ws.Range("D5") "215.58"
$ws.Range("E5").Value = "  +0.39%  "
Set-TextValue $ws.Range("D6") "0.5117"
$ws.Range("E6").Value = "  +0.38%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("E8").Value = "  -0.09%  "
Set-TextValue $ws.Range("D9") "0.06336"
$ws.Range("E9").Value = "  -0.38%  "
$ws.Range("E10").Value = "  -0.54%  "
$ws.Range("E11").Value = "  +0.46%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws.Range("D12") "1.639.23"
$ws.Range("E12").Value = "  -0.15%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue $ws.Range("D13") "4.241"
$ws.Range("E13").Value = "  -1.11%  "
Set-TextValue $ws.Range("D14") "1.853.22"
$ws.Range("E14").Value = "  -0.89%  "
Set-TextValue $ws.Range("D15") "0.5523"
$ws.Range("E15").Value = "  +1.44%  "
Set-TextValue $ws.Range("D16") "63.64"
$ws.Range("E16").Value = "  -0.99%  "
Set-TextValue $ws.Range("D17") "0.0₅7590"
$ws.Range("E17").Value = "  -1.78%  "
Set-TextValue $ws.Range("D18") "25.858.48"
$ws.Range("E18").Value = "  -0.64%  "
$ws.Range("E19").Value = "  -0.01%  "
Set-TextValue $ws.Range("D20") "194.57"
$ws.Range("E20").Value = "  -1.69%  "
$ws.Range("E21").Value = "  -0.24%  "
Set-TextValue $ws.Range("D22") "9.856"
$ws.Range("E22").Value = "  -0.74%  "
Set-TextValue $ws.Range("D23") "6.013"
$ws.Range("E23").Value = "  -0.45%  "
Set-TextValue $ws.Range("D24") "1.001"
$ws.Range("E24").Value = "  -0.18%  "
Set-TextValue $ws.Range("D25") "1.895"
$ws.Range("E25").Value = "  +1.62%  "
Set-TextValue $ws.Range("D26") "142.02"
$ws.Range("E26").Value = "  +0.80%  "
Set-TextValue $ws.Range("D27") "0.1253"
$ws.Range("E27").Value = "  +4.74%  "
Set-TextValue $ws.Range("D28") "6.756"
$ws.Range("E28").Value = "  -0.89%  "
Set-TextValue $ws.Range("D29") "15.56"
$ws.Range("E29").Value = "  +0.07%  "
Set-TextValue $ws.Range("D30") "1.238"
$ws.Range("E30").Value = "  +0.33%  "
Set-TextValue $ws.Range("D31") "0.04905"
$ws.Range("E31").Value = "  +0.96%  "
Set-TextValue $ws.Range("D32") "3.232"
$ws.Range("E32").Value = "  -0.71%  "
$ws.Range("E33").Value = "  +0.27%  "
Set-TextValue $ws.Range("D34") "1.545"
$ws.Range("E34").Value = "  +1.17%  "
Set-TextValue $ws.Range("D35") "2.371"
$ws.Range("E35").Value = "  +0.45%  "
Set-TextValue $ws.Range("D36") "0.8949"
$ws.Range("E36").Value = "  -0.52%  "
Set-TextValue $ws.Range("D37") "0.5518"
$ws.Range("E37").Value = "  +0.99%  "
Set-TextValue $ws.Range("D38") "2.533"
$ws.Range("E38").Value = "  -1.68%  "
Set-TextValue $ws.Range("D39") "1.115.37"
$ws.Range("E39").Value = "  -2.36%  "
Set-TextValue $ws.Range("D40") "0.01553"
$ws.Range("E40").Value = "  -0.52%  "
Set-TextValue $ws.Range("D41") "0.9997"
$ws.Range("E41").Value = "  -0.07%  "
Set-TextValue $ws.Range("D42") "5.567"
$ws.Range("E42").Value = "  +3.32%  "
Set-TextValue $ws.Range("D43") "0.7944"
$ws.Range("E43").Value = "  -2.00%  "
Set-TextValue $ws.Range("D44") "97.67"
$ws.Range("E44").Value = "  -1.72%  "
Set-TextValue $ws.Range("D45") "1.777.16"
$ws.Range("E45").Value = "  -0.13%  "
Set-TextValue $ws.Range("D46") "0.0₈116"
$ws.Range("E46").Value = "  -9.63%  "
Set-TextValue $ws.Range("D47") "0.4431"
$ws.Range("E47").Value = "  -1.99%  "
Set-TextValue $ws.Range("D48") "1.0000"
$ws.Range("E48").Value = "  +0.13%  "
$ws.Range("E49").Value = "  -0.38%  "
$ws.Range("E50").Value = "  +1.48%  "
Set-TextValue $ws.Range("D51") "7.578"
$ws.Range("E51").Value = "  +3.56%  "
